$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.622.17'
$ws.Range("E2").Value = '  -3.77%  '
$ws.Range("D3").Value = '2.542.69'
$ws.Range("E3").Value = '  -3.58%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '506.69'
$ws.Range("E5").Value = '  -4.26%  '
$ws.Range("D6").Value = '143.75'
$ws.Range("E6").Value = '  -7.33%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.564'
$ws.Range("E8").Value = '  -4.29%  '
$ws.Range("D9").Value = '2.549.28'
$ws.Range("E9").Value = '  -3.71%  '
$ws.Range("D10").Value = '6.10'
$ws.Range("E10").Value = '  -8.17%  '
$ws.Range("E12").Value = '  -5.52%  '
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("D14").Value = '2.987.80'
$ws.Range("E14").Value = '  -3.58%  '
$ws.Range("D15").Value = '58.585.15'
$ws.Range("E15").Value = '  -3.85%  '
$ws.Range("D16").Value = '20.73'
$ws.Range("E16").Value = '  -5.43%  '
$ws.Range("E17").Value = '  -5.65%  '
$ws.Range("D18").Value = '2.543.70'
$ws.Range("E18").Value = '  -3.57%  '
$ws.Range("E19").Value = '  -4.70%  '
$ws.Range("D20").Value = '340.54'
$ws.Range("E20").Value = '  -3.72%  '
$ws.Range("D21").Value = '10.11'
$ws.Range("E21").Value = '  -4.99%  '
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").Value = '5.97'
$ws.Range("E23").Value = '  -4.08%  '
$ws.Range("E24").Value = '  -1.37%  '
$ws.Range("D25").Value = '0.412'
$ws.Range("E25").Value = '  -4.24%  '
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  -0.24%  '
$ws.Range("E27").Value = '  -4.90%  '
$ws.Range("D28").Value = '2.653.10'
$ws.Range("E28").Value = '  -3.65%  '
$ws.Range("D29").Value = '0.0₃0788'
$ws.Range("E29").Value = '  -8.66%  '
$ws.Range("D30").Value = '6.97'
$ws.Range("E30").Value = '  -5.70%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = '149.68'
$ws.Range("E32").Value = '  -0.41%  '
$ws.Range("D33").Value = '5.84'
$ws.Range("E33").Value = '  -5.01%  '
$ws.Range("E34").Value = '  -4.80%  '
$ws.Range("D35").Value = '1.54'
$ws.Range("E35").Value = '  -5.59%  '
$ws.Range("D36").Value = '0.913'
$ws.Range("E36").Value = '  +2.13%  '
$ws.Range("E37").Value = '  -6.06%  '
$ws.Range("E38").Value = '  -7.36%  '
$ws.Range("D39").Value = '36.12'
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("E40").Value = '  -10.99%  '
$ws.Range("E41").Value = '  -6.81%  '
$ws.Range("D42").Value = '283.05'
$ws.Range("E42").Value = '  -7.88%  '
$ws.Range("E43").Value = '  -7.71%  '
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("D45").Value = '0.997'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '0.601'
$ws.Range("E46").Value = '  -6.18%  '
$ws.Range("D47").Value = '0.0532'
$ws.Range("E47").Value = '  -5.47%  '
$ws.Range("D48").Value = '18.70'
$ws.Range("E48").Value = '  -5.30%  '
$ws.Range("E49").Value = '  -0.50%  '
$ws.Range("D50").Value = '0.0227'
$ws.Range("E50").Value = '  -5.22%  '
$ws.Range("E51").Value = '  -8.62%  '
